$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate "2021-Q4" sheet as a formatting template, rename to "2022-Q1" ---
$srcQ4 = $wb.Worksheets.Item(4)
$srcQ4.Copy($null, $srcQ4)
$newWs = $wb.Worksheets.Item(5)
$newWs.Name = "2022-Q1"

# the template (2021-Q4) only has 15 data rows (rows 2-16); 2022-Q1 needs 16 (rows 2-17)
# extend the "index" column formatting (style carrying bold/border header look-alike for col A) down one more row
$newWs.Range("A16").Copy()
$newWs.Range("A17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Step 2: overwrite data rows with 2022-Q1 figures ---
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "'007689"
$newWs.Range("C2").Value = "国投瑞银新能源混合A"
$newWs.Range("D2").Value = "'78.74"
$newWs.Range("E2").Value = "'91.13"
$newWs.Range("F2").Value = "'4.65"
$newWs.Range("G2").Value = "'3.6614"
$newWs.Range("H2").Value = 8

$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "'012148"
$newWs.Range("C3").Value = "国投瑞银产业趋势混合型证券投资基金A"
$newWs.Range("D3").Value = "'45.53"
$newWs.Range("E3").Value = "'92.28"
$newWs.Range("F3").Value = "'5.49"
$newWs.Range("G3").Value = "'2.4996"
$newWs.Range("H3").Value = 5

$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "'006736"
$newWs.Range("C4").Value = "国投瑞银先进制造混合"
$newWs.Range("D4").Value = "'41.36"
$newWs.Range("E4").Value = "'92.33"
$newWs.Range("F4").Value = "'4.75"
$newWs.Range("G4").Value = "'1.9646"
$newWs.Range("H4").Value = 6

$newWs.Range("A5").Value = 3
$newWs.Range("B5").Value = "'007690"
$newWs.Range("C5").Value = "国投瑞银新能源混合C"
$newWs.Range("D5").Value = "'37.84"
$newWs.Range("E5").Value = "'91.13"
$newWs.Range("F5").Value = "'4.65"
$newWs.Range("G5").Value = "'1.7596"
$newWs.Range("H5").Value = 8

$newWs.Range("A6").Value = 4
$newWs.Range("B6").Value = "'011128"
$newWs.Range("C6").Value = "华安精致生活混合A"
$newWs.Range("D6").Value = "'33.30"
$newWs.Range("E6").Value = "'85.22"
$newWs.Range("F6").Value = "'3.61"
$newWs.Range("G6").Value = "'1.2021"
$newWs.Range("H6").Value = 2

$newWs.Range("A7").Value = 5
$newWs.Range("B7").Value = "'011251"
$newWs.Range("C7").Value = "华安聚嘉精选混合A"
$newWs.Range("D7").Value = "'31.29"
$newWs.Range("E7").Value = "'89.00"
$newWs.Range("F7").Value = "'3.54"
$newWs.Range("G7").Value = "'1.1077"
$newWs.Range("H7").Value = 1

$newWs.Range("A8").Value = 6
$newWs.Range("B8").Value = "'012149"
$newWs.Range("C8").Value = "国投瑞银产业趋势混合型证券投资基金C"
$newWs.Range("D8").Value = "'18.42"
$newWs.Range("E8").Value = "'92.28"
$newWs.Range("F8").Value = "'5.49"
$newWs.Range("G8").Value = "'1.0113"
$newWs.Range("H8").Value = 5

$newWs.Range("A9").Value = 7
$newWs.Range("B9").Value = "'007139"
$newWs.Range("C9").Value = "富国民裕进取沪港深成长精选混合"
$newWs.Range("D9").Value = "'12.79"
$newWs.Range("E9").Value = "'92.21"
$newWs.Range("F9").Value = "'6.92"
$newWs.Range("G9").Value = "'0.8851"
$newWs.Range("H9").Value = 3

$newWs.Range("A10").Value = 8
$newWs.Range("B10").Value = "'011252"
$newWs.Range("C10").Value = "华安聚嘉精选混合C"
$newWs.Range("D10").Value = "'14.89"
$newWs.Range("E10").Value = "'89.00"
$newWs.Range("F10").Value = "'3.54"
$newWs.Range("G10").Value = "'0.5271"
$newWs.Range("H10").Value = 1

$newWs.Range("A11").Value = 9
$newWs.Range("B11").Value = "'900008"
$newWs.Range("C11").Value = "中信证券稳健回报混合A"
$newWs.Range("D11").Value = "'6.36"
$newWs.Range("E11").Value = "'88.19"
$newWs.Range("F11").Value = "'7.80"
$newWs.Range("G11").Value = "'0.4961"
$newWs.Range("H11").Value = 1

$newWs.Range("A12").Value = 10
$newWs.Range("B12").Value = "'011129"
$newWs.Range("C12").Value = "华安精致生活混合C"
$newWs.Range("D12").Value = "'7.57"
$newWs.Range("E12").Value = "'85.22"
$newWs.Range("F12").Value = "'3.61"
$newWs.Range("G12").Value = "'0.2733"
$newWs.Range("H12").Value = 2

$newWs.Range("A13").Value = 11
$newWs.Range("B13").Value = "'900078"
$newWs.Range("C13").Value = "中信证券稳健回报混合C"
$newWs.Range("D13").Value = "'1.66"
$newWs.Range("E13").Value = "'88.19"
$newWs.Range("F13").Value = "'7.80"
$newWs.Range("G13").Value = "'0.1295"
$newWs.Range("H13").Value = 1

$newWs.Range("A14").Value = 12
$newWs.Range("B14").Value = "'008861"
$newWs.Range("C14").Value = "西部利得港股通新机遇灵活配置混合A"
$newWs.Range("D14").Value = "'0.37"
$newWs.Range("E14").Value = "'77.31"
$newWs.Range("F14").Value = "'3.76"
$newWs.Range("G14").Value = "'0.0139"
$newWs.Range("H14").Value = 6

$newWs.Range("A15").Value = 13
$newWs.Range("B15").Value = "'012315"
$newWs.Range("C15").Value = "创金合信港股通成长股票型发起式证券投资基金A"
$newWs.Range("D15").Value = "'0.19"
$newWs.Range("E15").Value = "'83.49"
$newWs.Range("F15").Value = "'6.61"
$newWs.Range("G15").Value = "'0.0126"
$newWs.Range("H15").Value = 6

$newWs.Range("A16").Value = 14
$newWs.Range("B16").Value = "'012316"
$newWs.Range("C16").Value = "创金合信港股通成长股票型发起式证券投资基金C"
$newWs.Range("D16").Value = "'0.10"
$newWs.Range("E16").Value = "'83.49"
$newWs.Range("F16").Value = "'6.61"
$newWs.Range("G16").Value = "'0.0066"
$newWs.Range("H16").Value = 6

$newWs.Range("A17").Value = 15
$newWs.Range("B17").Value = "'010093"
$newWs.Range("C17").Value = "西部利得港股通新机遇灵活配置混合C"
$newWs.Range("D17").Value = "'0.08"
$newWs.Range("E17").Value = "'77.31"
$newWs.Range("F17").Value = "'3.76"
$newWs.Range("G17").Value = "'0.0030"
$newWs.Range("H17").Value = 6

# --- Step 3: update the "总计" (totals) sheet: add a 2022-Q1 row on top, push existing rows down ---
$totalWs = $wb.Worksheets.Item($wb.Worksheets.Count)

# extend col-A index formatting one more row, then rewrite all 5 data rows top-to-bottom
$totalWs.Range("A5").Copy()
$totalWs.Range("A6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 16
$totalWs.Range("D2").Value = 15.55

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2021-Q4"
$totalWs.Range("C3").Value = 15
$totalWs.Range("D3").Value = 12.5

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2021-Q3"
$totalWs.Range("C4").Value = 46
$totalWs.Range("D4").Value = 36.59

$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q2"
$totalWs.Range("C5").Value = 2
$totalWs.Range("D5").Value = 1.81

$totalWs.Range("A6").Value = 4
$totalWs.Range("B6").Value = "2021-Q1"
$totalWs.Range("C6").Value = 2
$totalWs.Range("D6").Value = 1.35
